# Adds three label shapes to the Grays Harbor map figure:
#   - "Westport" textbox (flipped horizontally, no fill)
#   - a small accent-colored oval marker
#   - "Ocosta School" textbox (two runs: "Ocosta" + " School")
#
# PowerPoint's Shape.Left/Top/Width/Height (and the AddTextbox/AddShape
# position args) are expressed in points and stored internally as
# single-precision floats, while the OOXML stores EMUs (1 pt = 12700 EMU).
# A naive `emu / 12700.0` can truncate to one EMU below the target once
# the value is round-tripped through a 32-bit float, so EMU() nudges the
# point value up to the smallest float32 that still maps back to the
# exact target EMU.
function EMU($emu) {
    $base = $emu / 12700.0
    for ($k = 0; $k -le 4000; $k++) {
        $cand = $base + ($k * 0.0000002)
        $f = [single]$cand
        $back = [math]::Floor(([double]$f) * 12700.0 + 0.0000001)
        if ($back -eq $emu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Westport" label (TextBox 4, id=5) ---------------------------------
$westport = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$westport.TextFrame.WordWrap = -1
$westport.TextFrame.AutoSize = 1
$westport.TextFrame.TextRange.Text = "Westport"
$westport.Left = (EMU 104195)
$westport.Top = (EMU 2450068)
$westport.Width = (EMU 1102305)
$westport.Height = (EMU 369332)
$westport.HorizontalFlip = -1
$westport.Fill.Visible = 0

# --- Oval marker (Oval 9, id=10) ----------------------------------------
$marker = $s.Shapes.AddShape(9, 0, 0, 10, 10)
$marker.Left = (EMU 1282700)
$marker.Top = (EMU 2927767)
$marker.Width = (EMU 88900)
$marker.Height = (EMU 76299)
$marker.Fill.ForeColor.ObjectThemeColor = 6
$marker.Line.ForeColor.ObjectThemeColor = 6
$marker.TextFrame.VerticalAnchor = 3
$marker.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- "Ocosta School" label (TextBox 10, id=11) --------------------------
$ocosta = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$ocosta.TextFrame.WordWrap = -1
$ocosta.TextFrame.AutoSize = 1
$ocosta.TextFrame.TextRange.Text = "Ocosta"
[void]$ocosta.TextFrame.TextRange.InsertAfter(" School")
$ocosta.Left = (EMU 1371600)
$ocosta.Top = (EMU 2743101)
$ocosta.Width = (EMU 1778000)
$ocosta.Height = (EMU 369332)
$ocosta.Fill.Visible = 0

Write-Host "Added shapes:" $westport.Id $westport.Name "," $marker.Id $marker.Name "," $ocosta.Id $ocosta.Name
